$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") "63.708.33"
Set-TextValue $ws.Range("E2") "  +1.15%  "

# Row 3
Set-TextValue $ws.Range("D3") "3.283.89"
Set-TextValue $ws.Range("E3") "  +1.32%  "

# Row 4
Set-TextValue $ws.Range("D4") "1.00"
Set-TextValue $ws.Range("E4") "  +0.26%  "

# Row 5
Set-TextValue $ws.Range("D5") "535.76"
Set-TextValue $ws.Range("E5") "  +5.38%  "

# Row 6
Set-TextValue $ws.Range("D6") "174.89"
Set-TextValue $ws.Range("E6") "  +0.02%  "

# Row 7
Set-TextValue $ws.Range("D7") "0.598"
Set-TextValue $ws.Range("E7") "  +1.99%  "

# Row 8
Set-TextValue $ws.Range("D8") "3.283.73"
Set-TextValue $ws.Range("E8") "  +1.24%  "

# Row 9
Set-TextValue $ws.Range("D9") "1.00"
Set-TextValue $ws.Range("E9") "  -0.16%  "

# Row 10
Set-TextValue $ws.Range("D10") "0.612"
Set-TextValue $ws.Range("E10") "  +0.83%  "

# Row 11
Set-TextValue $ws.Range("D11") "53.97"
Set-TextValue $ws.Range("E11") "  -4.03%  "

# Row 12
Set-TextValue $ws.Range("E12") "  +6.03%  "

# Row 13
Set-TextValue $ws.Range("D13") "0.0000259"
Set-TextValue $ws.Range("E13") "  +3.37%  "

# Row 14
Set-TextValue $ws.Range("D14") "9.19"
Set-TextValue $ws.Range("E14") "  +2.18%  "

# Row 15
Set-TextValue $ws.Range("D15") "3.805.40"
Set-TextValue $ws.Range("E15") "  +2.06%  "

# Row 16
Set-TextValue $ws.Range("E16") "  +0.07%  "

# Row 17
Set-TextValue $ws.Range("D17") "3.271.40"
Set-TextValue $ws.Range("E17") "  +1.49%  "

# Row 18
Set-TextValue $ws.Range("E18") "  +1.89%  "

# Row 19
Set-TextValue $ws.Range("D19") "63.653.41"
Set-TextValue $ws.Range("E19") "  +1.39%  "

# Row 20
Set-TextValue $ws.Range("D20") "11.21"
Set-TextValue $ws.Range("E20") "  +4.37%  "

# Row 21
Set-TextValue $ws.Range("D21") "0.970"
Set-TextValue $ws.Range("E21") "  +4.35%  "

# Row 22
Set-TextValue $ws.Range("D22") "371.14"
Set-TextValue $ws.Range("E22") "  +1.68%  "

# Row 23
$ws.Range("B23").Value = "RenderToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Range("D23") "11.39"
Set-TextValue $ws.Range("E23") "  +4.33%  "

# Row 24
$ws.Range("B24").Value = "PancakeSwap"
$ws.Range("C24").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue $ws.Range("D24") "3.79"
Set-TextValue $ws.Range("E24") "  +5.84%  "

# Row 25
Set-TextValue $ws.Range("D25") "4.11"
Set-TextValue $ws.Range("E25") "  +9.56%  "

# Row 26
Set-TextValue $ws.Range("D26") "81.38"
Set-TextValue $ws.Range("E26") "  +3.52%  "

# Row 27
Set-TextValue $ws.Range("D27") "6.18"
Set-TextValue $ws.Range("E27") "  +4.26%  "

# Row 28
Set-TextValue $ws.Range("D28") "2.67"
Set-TextValue $ws.Range("E28") "  +3.18%  "

# Row 29
Set-TextValue $ws.Range("D29") "11.39"
Set-TextValue $ws.Range("E29") "  +2.43%  "

# Row 30
Set-TextValue $ws.Range("E30") "  +2.13%  "

# Row 31
Set-TextValue $ws.Range("D31") "28.78"
Set-TextValue $ws.Range("E31") "  +3.20%  "

# Row 32
Set-TextValue $ws.Range("D32") "641.14"
Set-TextValue $ws.Range("E32") "  +0.02%  "

# Row 33
Set-TextValue $ws.Range("D33") "6.50"
Set-TextValue $ws.Range("E33") "  -0.42%  "

# Row 34
Set-TextValue $ws.Range("D34") "11.32"
Set-TextValue $ws.Range("E34") "  +2.97%  "

# Row 35
Set-TextValue $ws.Range("D35") "0.108"
Set-TextValue $ws.Range("E35") "  +5.49%  "

# Row 36
Set-TextValue $ws.Range("D36") "57.17"
Set-TextValue $ws.Range("E36") "  -0.82%  "

# Row 37
Set-TextValue $ws.Range("E37") "  -0.11%  "

# Row 38
Set-TextValue $ws.Range("D38") "37.01"
Set-TextValue $ws.Range("E38") "  +4.78%  "

# Row 39
Set-TextValue $ws.Range("E39") "  +2.67%  "

# Row 40
Set-TextValue $ws.Range("D40") "0.0₃0743"
Set-TextValue $ws.Range("E40") "  +15.73%  "

# Row 41
Set-TextValue $ws.Range("D41") "1.00"
Set-TextValue $ws.Range("E41") "  +0.58%  "

# Row 42
Set-TextValue $ws.Range("E42") "  +2.84%  "

# Row 43
Set-TextValue $ws.Range("D43") "2.61"
Set-TextValue $ws.Range("E43") "  +9.82%  "

# Row 44
Set-TextValue $ws.Range("D44") "2.929.47"
Set-TextValue $ws.Range("E44") "  +2.83%  "

# Row 45
Set-TextValue $ws.Range("D45") "2.98"
Set-TextValue $ws.Range("E45") "  +7.93%  "

# Row 46
Set-TextValue $ws.Range("E46") "  +5.84%  "

# Row 47
Set-TextValue $ws.Range("D47") "0.0399"
Set-TextValue $ws.Range("E47") "  +6.50%  "

# Row 48
$ws.Range("B48").Value = "ApeXProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
Set-TextValue $ws.Range("D48") "3.10"
Set-TextValue $ws.Range("E48") "  +7.69%  "

# Row 49
$ws.Range("B49").Value = "ThetaToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
Set-TextValue $ws.Range("D49") "2.63"
Set-TextValue $ws.Range("E49") "  +0.28%  "

# Row 50
Set-TextValue $ws.Range("E50") "  +2.92%  "

# Row 51
Set-TextValue $ws.Range("D51") "135.16"
Set-TextValue $ws.Range("E51") "  +4.54%  "
